# Updates cryptos list with latest price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.051.49"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.572.28"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "2.577.19"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("E12").Value = "  +11.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.345"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").Value = "3.023.80"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "59.083.92"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.40%  "
$ws.Range("E17").Value = "  +3.63%  "
$ws.Range("D18").Value = "2.575.13"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "337.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  -3.26%  "
$ws.Range("E25").Value = "  +6.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "0.0₃0778"
$ws.Range("E29").Value = "  +2.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.871"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.871"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "293.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.04%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0976"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0536"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.90%  "
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("D51").Value = "1.946.62"
$ws.Range("E51").Value = "  -0.01%  "
